$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.885.33'
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").Value = '2.917.44'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = "'593.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.46%  '

# Row 6
$ws.Range("D6").Value = "'145.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").Value = "'0.507"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.81%  '

# Row 9
$ws.Range("D9").Value = "'6.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.05%  '

# Row 10
$ws.Range("E10").Value = '  +0.02%  '

# Row 11
$ws.Range("D11").Value = "'0.438"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.10%  '

# Row 13
$ws.Range("D13").Value = "'33.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.74%  '

# Row 14
$ws.Range("E14").Value = '  -0.03%  '

# Row 15
$ws.Range("D15").Value = '3.402.51'

# Row 16
$ws.Range("D16").Value = '60.885.56'
$ws.Range("E16").Value = '  +0.58%  '

# Row 17
$ws.Range("D17").Value = "'6.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.54%  '

# Row 18
$ws.Range("D18").Value = '2.920.73'
$ws.Range("E18").Value = '  +0.90%  '

# Row 19
$ws.Range("D19").Value = "'430.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.45%  '

# Row 20
$ws.Range("D20").Value = "'13.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.89%  '

# Row 21
$ws.Range("E21").Value = '  +1.51%  '

# Row 22
$ws.Range("E22").Value = '  -0.33%  '

# Row 23
$ws.Range("D23").Value = "'81.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.59%  '

# Row 24
$ws.Range("D24").Value = "'10.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '

# Row 25
$ws.Range("D25").Value = "'2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.56%  '

# Row 26
$ws.Range("D26").Value = "'11.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.58%  '

# Row 27
$ws.Range("E27").Value = '  -0.01%  '

# Row 28
$ws.Range("D28").Value = "'2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.38%  '

# Row 29
$ws.Range("E29").Value = '  -0.01%  '

# Row 30
$ws.Range("E30").Value = '  -0.02%  '

# Row 31
$ws.Range("E31").Value = '  -3.53%  '

# Row 32
$ws.Range("D32").Value = "'26.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.09%  '

# Row 33
$ws.Range("E33").Value = '  +0.26%  '

# Row 34
$ws.Range("E34").Value = '  +2.04%  '

# Row 35
$ws.Range("E35").Value = '  +0.33%  '

# Row 36
$ws.Range("E36").Value = '  -1.04%  '

# Row 37
$ws.Range("E37").Value = '  +2.83%  '

# Row 38
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = "'1.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.76%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = "'0.122"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.34%  '

# Row 40
$ws.Range("E40").Value = '  -1.93%  '

# Row 41
$ws.Range("E41").Value = '  -2.02%  '

# Row 42
$ws.Range("D42").Value = "'40.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.94%  '

# Row 43
$ws.Range("D43").Value = "'373.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.92%  '

# Row 44
$ws.Range("E44").Value = '  -0.18%  '

# Row 45
$ws.Range("D45").Value = '2.696.39'
$ws.Range("E45").Value = '  +1.73%  '

# Row 46
$ws.Range("D46").Value = "'130.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '

# Row 48
$ws.Range("D48").Value = "'23.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.62%  '

# Row 49
$ws.Range("E49").Value = '  -0.25%  '

# Row 50
$ws.Range("E50").Value = '  -3.70%  '

# Row 51
$ws.Range("E51").Value = '  +2.06%  '
